$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells are stored as text in the source sheet (prices with
# thousands separators typed as literal dots). A handful of the new
# values are "plain" decimals that Excel would otherwise auto-convert
# to a number and silently drop a trailing zero from (e.g. 0.840 ->
# 0.84), so force those specific cells to Text before writing, then
# drop back to the default style so no stray formatting is left behind.

$ws.Range("D2").Value = "43.536.95"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").Value = "2.284.45"
$ws.Range("E3").Value = "  -1.29%  "

$ws.Range("D5").Value = "95.44"
$ws.Range("E5").Value = "  +1.24%  "

$ws.Range("D6").Value = "267.91"
$ws.Range("E6").Value = "  -1.13%  "

$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  -1.10%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("E9").Value = "  -2.37%  "

$ws.Range("D10").Value = "45.28"
$ws.Range("E10").Value = "  +0.74%  "

$ws.Range("D11").Value = "0.0932"
$ws.Range("E11").Value = "  -0.86%  "

$ws.Range("D12").Value = "7.91"
$ws.Range("E12").Value = "  -2.81%  "

$ws.Range("E13").Value = "  +0.74%  "

$ws.Range("D14").Value = "2.626.35"
$ws.Range("E14").Value = "  -1.27%  "

$ws.Range("D15").Value = "15.27"
$ws.Range("E15").Value = "  -0.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.840"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.97%  "

$ws.Range("D17").Value = "2.285.04"
$ws.Range("E17").Value = "  -1.08%  "

$ws.Range("D18").Value = "43.510.17"
$ws.Range("E18").Value = "  -0.62%  "

$ws.Range("E19").Value = "  +2.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.47%  "

$ws.Range("D21").Value = "71.99"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("D22").Value = "2.58"
$ws.Range("E22").Value = "  +12.56%  "

$ws.Range("D23").Value = "232.13"
$ws.Range("E23").Value = "  -2.94%  "

$ws.Range("D24").Value = "9.11"
$ws.Range("E24").Value = "  -6.24%  "

$ws.Range("E25").Value = "  -0.14%  "

$ws.Range("D26").Value = "2.55"
$ws.Range("E26").Value = "  +1.56%  "

$ws.Range("D27").Value = "11.18"
$ws.Range("E27").Value = "  -1.63%  "

$ws.Range("E28").Value = "  +2.19%  "

$ws.Range("D29").Value = "40.03"
$ws.Range("E29").Value = "  +2.25%  "

$ws.Range("E30").Value = "  -2.11%  "

$ws.Range("D31").Value = "174.99"
$ws.Range("E31").Value = "  +1.61%  "

$ws.Range("D32").Value = "21.76"
$ws.Range("E32").Value = "  -3.66%  "

$ws.Range("D33").Value = "0.0893"
$ws.Range("E33").Value = "  -1.07%  "

$ws.Range("E34").Value = "  -4.14%  "

$ws.Range("E35").Value = "  -0.79%  "

$ws.Range("E36").Value = "  -3.54%  "

$ws.Range("E37").Value = "  -1.98%  "

$ws.Range("E38").Value = "  -3.06%  "

$ws.Range("E39").Value = "  -4.06%  "

$ws.Range("D40").Value = "0.238"
$ws.Range("E40").Value = "  +0.64%  "

$ws.Range("E41").Value = "  +0.11%  "

$ws.Range("D42").Value = "12.26"
$ws.Range("E42").Value = "  +0.42%  "

$ws.Range("D43").Value = "65.54"
$ws.Range("E43").Value = "  +5.76%  "

$ws.Range("D44").Value = "1.34"
$ws.Range("E44").Value = "  +1.18%  "

$ws.Range("D45").Value = "8.75"
$ws.Range("E45").Value = "  -2.73%  "

$ws.Range("E46").Value = "  -5.73%  "

$ws.Range("D47").Value = "0.102"
$ws.Range("E47").Value = "  -0.71%  "

$ws.Range("E48").Value = "  -1.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "95.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.68%  "

$ws.Range("E50").Value = "  +7.88%  "

$ws.Range("D51").Value = "0.425"
$ws.Range("E51").Value = "  -0.72%  "
